$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (Exhibition) ---
$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("B3").Value = "'2024-08-24"
$ws1.Range("C3").Value = "苏州·幻想物语次元嘉年华（免费展）"
$ws1.Range("D3").Value = "相城大道1168号 天虹购物中心(相城店)"
$ws1.Range("E3").Value = "2024.08.24 14:00-08.25 17:00"
$ws1.Range("F3").Value = 1072
$ws1.Range("G3").Value = 40
$ws1.Range("H3").Value = "https://show.bilibili.com/platform/detail.html?id=90768"
$ws1.Range("I3").Value = "//i2.hdslb.com/bfs/openplatform/202408/OZRx7O051723788701169.jpeg"
$ws1.Range("B4").Value = "'2024-08-24"
$ws1.Range("C4").Value = "苏州·排球少年only-茶歇"
$ws1.Range("D4").Value = "德必姑苏WE国际文化艺术中心6-102室渔郎桥浜路16号 星渡咖啡"
$ws1.Range("E4").Value = "2024.08.24 10:00-08.25 19:00"
$ws1.Range("F4").Value = 407
$ws1.Range("G4").Value = 50
$ws1.Range("H4").Value = "https://show.bilibili.com/platform/detail.html?id=88689"
$ws1.Range("I4").Value = "//i0.hdslb.com/bfs/openplatform/202407/dX8i0duL1720087129283.jpeg"
$ws1.Range("B5").Value = "'2024-08-24"
$ws1.Range("C5").Value = "苏州·星梦X动漫游戏展（免费展）"
$ws1.Range("D5").Value = "劳动路725号 首开龙湖苏州胥江天街"
$ws1.Range("E5").Value = "2024.08.24 10:00-08.24 17:00"
$ws1.Range("F5").Value = 242
$ws1.Range("G5").Value = 20
$ws1.Range("H5").Value = "https://show.bilibili.com/platform/detail.html?id=88825"
$ws1.Range("I5").Value = "//i1.hdslb.com/bfs/openplatform/202407/WIfUXUJ81720411118529.jpeg"
$ws1.Range("B6").Value = "'2024-08-24"
$ws1.Range("C6").Value = "苏州·赛马娘ONLY"
$ws1.Range("D6").Value = "东苑路115-11号 苏苑饭店"
$ws1.Range("E6").Value = "2024.08.24 10:00-08.24 16:00"
$ws1.Range("F6").Value = 139
$ws1.Range("G6").Value = 60
$ws1.Range("H6").Value = "https://show.bilibili.com/platform/detail.html?id=87619"
$ws1.Range("I6").Value = "//i2.hdslb.com/bfs/openplatform/202406/2seg6U5l1718215220516.jpeg"
$ws1.Range("B7").Value = "'2024-09-15"
$ws1.Range("C7").Value = "苏州·Good jump ACG中秋嘉年华动漫国潮文化节"
$ws1.Range("D7").Value = "金山南路288号木渎影视城F2 苏州广电国际会展中心"
$ws1.Range("E7").Value = "2024.09.15 10:00-09.16 17:00"
$ws1.Range("F7").Value = 11976
$ws1.Range("G7").Value = 60
$ws1.Range("H7").Value = "https://show.bilibili.com/platform/detail.html?id=87120"
$ws1.Range("I7").Value = "//i0.hdslb.com/bfs/openplatform/202407/yw21E7Vn1721701909995.jpeg"
$ws1.Range("B8").Value = "'2024-09-15"
$ws1.Range("C8").Value = "苏州·第二届-AME动漫嘉年华（免费展）"
$ws1.Range("D8").Value = "东吴南路179号 龙湖苏州东吴天街"
$ws1.Range("E8").Value = "2024.09.15 12:00-09.15 16:30"
$ws1.Range("F8").Value = 41
$ws1.Range("G8").Value = 39.9
$ws1.Range("H8").Value = "https://show.bilibili.com/platform/detail.html?id=90288"
$ws1.Range("I8").Value = "//i0.hdslb.com/bfs/openplatform/202408/e3uZDian1722619198829.png"
$ws1.Range("B9").Value = "'2024-10-01"
$ws1.Range("C9").Value = "【大会员提前抢】苏州·OCG国潮动漫游戏嘉年华CV杨天翔内场"
$ws1.Range("D9").Value = "苏州大道东688号 苏州国际博览中心"
$ws1.Range("E9").Value = "2024.10.01 09:00-10.01 17:00"
$ws1.Range("F9").Value = 69
$ws1.Range("G9").Value = 258
$ws1.Range("H9").Value = "https://show.bilibili.com/platform/detail.html?id=90770"
$ws1.Range("I9").Value = "//i1.hdslb.com/bfs/openplatform/202408/reVAMzAd1723703136204.jpeg"
$ws1.Range("B10").Value = "'2024-10-01"
$ws1.Range("C10").Value = "苏州·I COME ACG动漫品牌博览会"
$ws1.Range("D10").Value = "金山南路288号木渎影视城F2 苏州广电国际会展中心"
$ws1.Range("E10").Value = "2024.10.01 10:00-10.03 17:00"
$ws1.Range("F10").Value = 11694
$ws1.Range("G10").Value = 60
$ws1.Range("H10").Value = "https://show.bilibili.com/platform/detail.html?id=87118"
$ws1.Range("I10").Value = "//i2.hdslb.com/bfs/openplatform/202406/SyK3KnNb1717829071115.jpeg"
$ws1.Range("B11").Value = "'2024-10-01"
$ws1.Range("C11").Value = "苏州·理想乡动漫游戏展-两馆全开+三馆间通道"
$ws1.Range("D11").Value = "花桥经济开发区绿地大道1598号 花桥国际博览中心"
$ws1.Range("E11").Value = "2024.10.01 10:00-10.03 17:00"
$ws1.Range("F11").Value = 4731
$ws1.Range("G11").Value = 39
$ws1.Range("H11").Value = "https://show.bilibili.com/platform/detail.html?id=83821"
$ws1.Range("I11").Value = "//i0.hdslb.com/bfs/openplatform/202404/OMtuTTFY1711958198579.jpeg"
$ws1.Range("B12").Value = "'2024-10-01"
$ws1.Range("C12").Value = "苏州·第四届-OCG国朝动漫游戏嘉年华"
$ws1.Range("D12").Value = "苏州大道东688号 苏州国际博览中心"
$ws1.Range("E12").Value = "2024.10.01 09:00-10.02 17:00"
$ws1.Range("F12").Value = 523
$ws1.Range("G12").Value = 70
$ws1.Range("H12").Value = "https://show.bilibili.com/platform/detail.html?id=89473"
$ws1.Range("I12").Value = "//i1.hdslb.com/bfs/openplatform/202407/q1QT4jWI1720809490212.jpeg"
$ws1.Range("B13").Value = "'2024-10-02"
$ws1.Range("C13").Value = "常熟·CDW·动漫展06"
$ws1.Range("D13").Value = "凯文路8号 常熟市体育中心-东北门"
$ws1.Range("E13").Value = "2024.10.02 09:00-10.03 17:30"
$ws1.Range("F13").Value = 63
$ws1.Range("G13").Value = 60
$ws1.Range("H13").Value = "https://show.bilibili.com/platform/detail.html?id=90736"
$ws1.Range("I13").Value = "//i1.hdslb.com/bfs/openplatform/202408/quQeX6tm1722926244359.jpeg"
$ws1.Range("B14").Value = "'2024-10-02"
$ws1.Range("C14").Value = "苏州·明日方舟ONLY#2024~佑桑柔"
$ws1.Range("D14").Value = "城际路21号 苏州汇融广场假日酒店"
$ws1.Range("E14").Value = "2024.10.02 10:00-10.02 17:00"
$ws1.Range("F14").Value = 418
$ws1.Range("G14").Value = 75
$ws1.Range("H14").Value = "https://show.bilibili.com/platform/detail.html?id=84046"
$ws1.Range("I14").Value = "//i2.hdslb.com/bfs/openplatform/202405/0VhIRprD1716344515303.jpeg"
$ws1.Range("B15").Value = "'2024-10-02"
$ws1.Range("C15").Value = "苏州·第二届百合Only同人展交流"
$ws1.Range("D15").Value = "三香路488号 苏州金陵雅都大酒店"
$ws1.Range("E15").Value = "2024.10.02 10:00-10.02 16:00"
$ws1.Range("F15").Value = 83
$ws1.Range("G15").Value = 46
$ws1.Range("H15").Value = "https://show.bilibili.com/platform/detail.html?id=89946"
$ws1.Range("I15").Value = "//i0.hdslb.com/bfs/openplatform/202407/nwznBIxG1721628287653.jpeg"
$ws1.Range("B16").Value = "'2024-10-04"
$ws1.Range("C16").Value = "常熟·cc动漫游戏嘉年华"
$ws1.Range("D16").Value = "开元大道1号 常熟国际博览中心"
$ws1.Range("E16").Value = "2024.10.04 09:00-10.05 17:00"
$ws1.Range("F16").Value = 921
$ws1.Range("G16").Value = 60
$ws1.Range("H16").Value = "https://show.bilibili.com/platform/detail.html?id=90292"
$ws1.Range("I16").Value = "//i2.hdslb.com/bfs/openplatform/202407/yCNXedrA1722404050722.jpeg"
$ws1.Range("B17").Value = "'2024-10-19"
$ws1.Range("C17").Value = "苏州·代号鸢only茶话会-星渡咖啡"
$ws1.Range("D17").Value = "德必姑苏WE国际文化艺术中心6-102室渔郎桥浜路16号 星渡咖啡"
$ws1.Range("E17").Value = "2024.10.19 10:00-10.20 19:00"
$ws1.Range("F17").Value = 351
$ws1.Range("G17").Value = 50
$ws1.Range("H17").Value = "https://show.bilibili.com/platform/detail.html?id=87685"
$ws1.Range("I17").Value = "//i1.hdslb.com/bfs/openplatform/202406/eyHRVQuv1718780132754.jpeg"
$ws1.Range("B18").Value = "'2024-10-26"
$ws1.Range("C18").Value = "苏州·第三届华盟国漫次元嘉年华"
$ws1.Range("D18").Value = "清禾路886号 苏州聚橙尹山湖大剧院"
$ws1.Range("E18").Value = "2024.10.26 10:00-10.27 17:00"
$ws1.Range("F18").Value = 157
$ws1.Range("G18").Value = 58
$ws1.Range("H18").Value = "https://show.bilibili.com/platform/detail.html?id=85767"
$ws1.Range("I18").Value = "//i1.hdslb.com/bfs/openplatform/202405/CqSYBZhQ1715846719965.jpeg"
$ws1.Range("B19").Value = "'2024-11-16"
$ws1.Range("C19").Value = "张家港·META萌圆饿了"
$ws1.Range("D19").Value = "杨舍镇人民中路42号 张家港国贸酒店"
$ws1.Range("E19").Value = "2024.11.16 10:00-11.16 17:00"
$ws1.Range("F19").Value = 51
$ws1.Range("G19").Value = 40
$ws1.Range("H19").Value = "https://show.bilibili.com/platform/detail.html?id=90745"
$ws1.Range("I19").Value = "//i2.hdslb.com/bfs/openplatform/202408/jB7b2kZ11723621730632.jpeg"
$ws1.Range("B20").Value = "'2024-12-27"
$ws1.Range("C20").Value = "苏州·星部落动漫嘉年华"
$ws1.Range("D20").Value = "花桥经济开发区绿地大道1598号 花桥国际博览中心"
$ws1.Range("E20").Value = "2024.12.27 09:00-12.28 16:00"
$ws1.Range("F20").Value = 5217
$ws1.Range("G20").Value = "不可售"
$ws1.Range("H20").Value = "https://show.bilibili.com/platform/detail.html?id=84858"
$ws1.Range("I20").Value = "//i0.hdslb.com/bfs/openplatform/202404/UI5EFZTT1713685680462.jpeg"
$ws1.Rows.Item(21).Delete()

# --- Sheet "全部类型" (All types) ---
$ws4 = $wb.Worksheets.Item(4)
$ws4.Range("B3").Value = "'2024-08-24"
$ws4.Range("C3").Value = "苏州·幻想物语次元嘉年华（免费展）"
$ws4.Range("D3").Value = "相城大道1168号 天虹购物中心(相城店)"
$ws4.Range("E3").Value = "2024.08.24 14:00-08.25 17:00"
$ws4.Range("F3").Value = 1072
$ws4.Range("G3").Value = 40
$ws4.Range("H3").Value = "https://show.bilibili.com/platform/detail.html?id=90768"
$ws4.Range("I3").Value = "//i2.hdslb.com/bfs/openplatform/202408/OZRx7O051723788701169.jpeg"
$ws4.Range("B4").Value = "'2024-08-24"
$ws4.Range("C4").Value = "苏州·排球少年only-茶歇"
$ws4.Range("D4").Value = "德必姑苏WE国际文化艺术中心6-102室渔郎桥浜路16号 星渡咖啡"
$ws4.Range("E4").Value = "2024.08.24 10:00-08.25 19:00"
$ws4.Range("F4").Value = 407
$ws4.Range("G4").Value = 50
$ws4.Range("H4").Value = "https://show.bilibili.com/platform/detail.html?id=88689"
$ws4.Range("I4").Value = "//i0.hdslb.com/bfs/openplatform/202407/dX8i0duL1720087129283.jpeg"
$ws4.Range("B5").Value = "'2024-08-24"
$ws4.Range("C5").Value = "苏州·星梦X动漫游戏展（免费展）"
$ws4.Range("D5").Value = "劳动路725号 首开龙湖苏州胥江天街"
$ws4.Range("E5").Value = "2024.08.24 10:00-08.24 17:00"
$ws4.Range("F5").Value = 242
$ws4.Range("G5").Value = 20
$ws4.Range("H5").Value = "https://show.bilibili.com/platform/detail.html?id=88825"
$ws4.Range("I5").Value = "//i1.hdslb.com/bfs/openplatform/202407/WIfUXUJ81720411118529.jpeg"
$ws4.Range("B6").Value = "'2024-08-24"
$ws4.Range("C6").Value = "苏州·赛马娘ONLY"
$ws4.Range("D6").Value = "东苑路115-11号 苏苑饭店"
$ws4.Range("E6").Value = "2024.08.24 10:00-08.24 16:00"
$ws4.Range("F6").Value = 139
$ws4.Range("G6").Value = 60
$ws4.Range("H6").Value = "https://show.bilibili.com/platform/detail.html?id=87619"
$ws4.Range("I6").Value = "//i2.hdslb.com/bfs/openplatform/202406/2seg6U5l1718215220516.jpeg"
$ws4.Range("B7").Value = "'2024-08-25"
$ws4.Range("C7").Value = "苏州·奇迹の闪耀 「UP!」巡回动漫演唱会"
$ws4.Range("D7").Value = "富乐路1号 阳澄文体中心"
$ws4.Range("E7").Value = "2024.08.25 19:30-08.25 21:30"
$ws4.Range("F7").Value = 4
$ws4.Range("G7").Value = 180
$ws4.Range("H7").Value = "https://show.bilibili.com/platform/detail.html?id=90249"
$ws4.Range("I7").Value = "//i0.hdslb.com/bfs/openplatform/202408/izTSQsGS1722582891356.jpeg"
$ws4.Range("B8").Value = "'2024-09-06"
$ws4.Range("C8").Value = "苏州·【明星版】吴琼主演经典黄梅戏《女驸马》"
$ws4.Range("D8").Value = "星湖街555号 苏州独墅湖影剧院"
$ws4.Range("E8").Value = "2024.09.06 19:30-09.07 21:30"
$ws4.Range("F8").Value = 0
$ws4.Range("G8").Value = 140
$ws4.Range("H8").Value = "https://show.bilibili.com/platform/detail.html?id=90156"
$ws4.Range("I8").Value = "//i0.hdslb.com/bfs/openplatform/202407/ZzBQjFLb1722413701814.jpeg"
$ws4.Range("B9").Value = "'2024-09-15"
$ws4.Range("C9").Value = "苏州·Good jump ACG中秋嘉年华动漫国潮文化节"
$ws4.Range("D9").Value = "金山南路288号木渎影视城F2 苏州广电国际会展中心"
$ws4.Range("E9").Value = "2024.09.15 10:00-09.16 17:00"
$ws4.Range("F9").Value = 11976
$ws4.Range("G9").Value = 60
$ws4.Range("H9").Value = "https://show.bilibili.com/platform/detail.html?id=87120"
$ws4.Range("I9").Value = "//i0.hdslb.com/bfs/openplatform/202407/yw21E7Vn1721701909995.jpeg"
$ws4.Range("B10").Value = "'2024-09-15"
$ws4.Range("C10").Value = "苏州·第二届-AME动漫嘉年华（免费展）"
$ws4.Range("D10").Value = "东吴南路179号 龙湖苏州东吴天街"
$ws4.Range("E10").Value = "2024.09.15 12:00-09.15 16:30"
$ws4.Range("F10").Value = 41
$ws4.Range("G10").Value = 39.9
$ws4.Range("H10").Value = "https://show.bilibili.com/platform/detail.html?id=90288"
$ws4.Range("I10").Value = "//i0.hdslb.com/bfs/openplatform/202408/e3uZDian1722619198829.png"
$ws4.Range("B11").Value = "'2024-10-01"
$ws4.Range("C11").Value = "【大会员提前抢】苏州·OCG国潮动漫游戏嘉年华CV杨天翔内场"
$ws4.Range("D11").Value = "苏州大道东688号 苏州国际博览中心"
$ws4.Range("E11").Value = "2024.10.01 09:00-10.01 17:00"
$ws4.Range("F11").Value = 69
$ws4.Range("G11").Value = 258
$ws4.Range("H11").Value = "https://show.bilibili.com/platform/detail.html?id=90770"
$ws4.Range("I11").Value = "//i1.hdslb.com/bfs/openplatform/202408/reVAMzAd1723703136204.jpeg"
$ws4.Range("B12").Value = "'2024-10-01"
$ws4.Range("C12").Value = "苏州·I COME ACG动漫品牌博览会"
$ws4.Range("D12").Value = "金山南路288号木渎影视城F2 苏州广电国际会展中心"
$ws4.Range("E12").Value = "2024.10.01 10:00-10.03 17:00"
$ws4.Range("F12").Value = 11694
$ws4.Range("G12").Value = 60
$ws4.Range("H12").Value = "https://show.bilibili.com/platform/detail.html?id=87118"
$ws4.Range("I12").Value = "//i2.hdslb.com/bfs/openplatform/202406/SyK3KnNb1717829071115.jpeg"
$ws4.Range("B13").Value = "'2024-10-01"
$ws4.Range("C13").Value = "苏州·理想乡动漫游戏展-两馆全开+三馆间通道"
$ws4.Range("D13").Value = "花桥经济开发区绿地大道1598号 花桥国际博览中心"
$ws4.Range("E13").Value = "2024.10.01 10:00-10.03 17:00"
$ws4.Range("F13").Value = 4731
$ws4.Range("G13").Value = 39
$ws4.Range("H13").Value = "https://show.bilibili.com/platform/detail.html?id=83821"
$ws4.Range("I13").Value = "//i0.hdslb.com/bfs/openplatform/202404/OMtuTTFY1711958198579.jpeg"
$ws4.Range("B14").Value = "'2024-10-01"
$ws4.Range("C14").Value = "苏州·第四届-OCG国朝动漫游戏嘉年华"
$ws4.Range("D14").Value = "苏州大道东688号 苏州国际博览中心"
$ws4.Range("E14").Value = "2024.10.01 09:00-10.02 17:00"
$ws4.Range("F14").Value = 523
$ws4.Range("G14").Value = 70
$ws4.Range("H14").Value = "https://show.bilibili.com/platform/detail.html?id=89473"
$ws4.Range("I14").Value = "//i1.hdslb.com/bfs/openplatform/202407/q1QT4jWI1720809490212.jpeg"
$ws4.Range("B15").Value = "'2024-10-02"
$ws4.Range("C15").Value = "常熟·CDW·动漫展06"
$ws4.Range("D15").Value = "凯文路8号 常熟市体育中心-东北门"
$ws4.Range("E15").Value = "2024.10.02 09:00-10.03 17:30"
$ws4.Range("F15").Value = 63
$ws4.Range("G15").Value = 60
$ws4.Range("H15").Value = "https://show.bilibili.com/platform/detail.html?id=90736"
$ws4.Range("I15").Value = "//i1.hdslb.com/bfs/openplatform/202408/quQeX6tm1722926244359.jpeg"
$ws4.Range("B16").Value = "'2024-10-02"
$ws4.Range("C16").Value = "苏州·明日方舟ONLY#2024~佑桑柔"
$ws4.Range("D16").Value = "城际路21号 苏州汇融广场假日酒店"
$ws4.Range("E16").Value = "2024.10.02 10:00-10.02 17:00"
$ws4.Range("F16").Value = 418
$ws4.Range("G16").Value = 75
$ws4.Range("H16").Value = "https://show.bilibili.com/platform/detail.html?id=84046"
$ws4.Range("I16").Value = "//i2.hdslb.com/bfs/openplatform/202405/0VhIRprD1716344515303.jpeg"
$ws4.Range("B17").Value = "'2024-10-02"
$ws4.Range("C17").Value = "苏州·第二届百合Only同人展交流"
$ws4.Range("D17").Value = "三香路488号 苏州金陵雅都大酒店"
$ws4.Range("E17").Value = "2024.10.02 10:00-10.02 16:00"
$ws4.Range("F17").Value = 83
$ws4.Range("G17").Value = 46
$ws4.Range("H17").Value = "https://show.bilibili.com/platform/detail.html?id=89946"
$ws4.Range("I17").Value = "//i0.hdslb.com/bfs/openplatform/202407/nwznBIxG1721628287653.jpeg"
$ws4.Range("B18").Value = "'2024-10-04"
$ws4.Range("C18").Value = "常熟·cc动漫游戏嘉年华"
$ws4.Range("D18").Value = "开元大道1号 常熟国际博览中心"
$ws4.Range("E18").Value = "2024.10.04 09:00-10.05 17:00"
$ws4.Range("F18").Value = 921
$ws4.Range("G18").Value = 60
$ws4.Range("H18").Value = "https://show.bilibili.com/platform/detail.html?id=90292"
$ws4.Range("I18").Value = "//i2.hdslb.com/bfs/openplatform/202407/yCNXedrA1722404050722.jpeg"
$ws4.Range("B19").Value = "'2024-10-19"
$ws4.Range("C19").Value = "苏州·代号鸢only茶话会-星渡咖啡"
$ws4.Range("D19").Value = "德必姑苏WE国际文化艺术中心6-102室渔郎桥浜路16号 星渡咖啡"
$ws4.Range("E19").Value = "2024.10.19 10:00-10.20 19:00"
$ws4.Range("F19").Value = 351
$ws4.Range("G19").Value = 50
$ws4.Range("H19").Value = "https://show.bilibili.com/platform/detail.html?id=87685"
$ws4.Range("I19").Value = "//i1.hdslb.com/bfs/openplatform/202406/eyHRVQuv1718780132754.jpeg"
$ws4.Range("B20").Value = "'2024-10-26"
$ws4.Range("C20").Value = "苏州·第三届华盟国漫次元嘉年华"
$ws4.Range("D20").Value = "清禾路886号 苏州聚橙尹山湖大剧院"
$ws4.Range("E20").Value = "2024.10.26 10:00-10.27 17:00"
$ws4.Range("F20").Value = 157
$ws4.Range("G20").Value = 58
$ws4.Range("H20").Value = "https://show.bilibili.com/platform/detail.html?id=85767"
$ws4.Range("I20").Value = "//i1.hdslb.com/bfs/openplatform/202405/CqSYBZhQ1715846719965.jpeg"
$ws4.Range("B21").Value = "'2024-11-16"
$ws4.Range("C21").Value = "张家港·META萌圆饿了"
$ws4.Range("D21").Value = "杨舍镇人民中路42号 张家港国贸酒店"
$ws4.Range("E21").Value = "2024.11.16 10:00-11.16 17:00"
$ws4.Range("F21").Value = 51
$ws4.Range("G21").Value = 40
$ws4.Range("H21").Value = "https://show.bilibili.com/platform/detail.html?id=90745"
$ws4.Range("I21").Value = "//i2.hdslb.com/bfs/openplatform/202408/jB7b2kZ11723621730632.jpeg"
$ws4.Range("B22").Value = "'2024-12-27"
$ws4.Range("C22").Value = "苏州·星部落动漫嘉年华"
$ws4.Range("D22").Value = "花桥经济开发区绿地大道1598号 花桥国际博览中心"
$ws4.Range("E22").Value = "2024.12.27 09:00-12.28 16:00"
$ws4.Range("F22").Value = 5217
$ws4.Range("G22").Value = "不可售"
$ws4.Range("H22").Value = "https://show.bilibili.com/platform/detail.html?id=84858"
$ws4.Range("I22").Value = "//i0.hdslb.com/bfs/openplatform/202404/UI5EFZTT1713685680462.jpeg"
$ws4.Rows.Item(23).Delete()
